$d = $word.ActiveDocument

$para = $d.Paragraphs.Item(1)
$full = $d.Content
$textRange = $d.Range($para.Range.Start, $full.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ФОРАШРУАПЗШЦртагуцращШГПРАКЩ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>СТШГАЩШГЦпадфы</w:t></w:r><w:r><w:t>ВАСФ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ТЗАшгфыунзашфцрАЗТуцгнкс8Йнцвзщш НИШЗЩАЦЙнважзщгтщцфшТВГСсъ</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>0</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/><w:t>ухтйВфщыгвстзшГВж0</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$textRange.InsertXML($xml)
